# Update the "Jogos da Semana" FlashScore workbook with refreshed odds values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 updates
$ws.Range("G3").Value = 3
$ws.Range("I3").Value = 2.3
$ws.Range("J3").Value = 3.6
$ws.Range("L3").Value = 3
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.9
$ws.Range("U3").Value = 1.7
$ws.Range("V3").Value = 2.05
$ws.Range("AB3").Value = 29
$ws.Range("AC3").Value = 11
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 12
$ws.Range("AJ3").Value = 9.5
$ws.Range("AP3").Value = 23
$ws.Range("AR3").Value = 67
$ws.Range("AU3").Value = 7.5

# Row 7 updates
$ws.Range("I7").Value = 5.5
$ws.Range("Q7").Value = 2.6
$ws.Range("R7").Value = 1.48
$ws.Range("W7").Value = 4.75
$ws.Range("AA7").Value = 19
$ws.Range("AE7").Value = 26
$ws.Range("AM7").Value = 67
$ws.Range("AP7").Value = 29
